$d = $word.ActiveDocument
$Q = [char]34

function Set-ParagraphText($matchSubstring, $newText) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like $matchSubstring) {
            $rng = $p.Range
            $rng.End = $rng.End - 1
            $rng.Text = $newText
            return $i
        }
    }
    return -1
}

function Insert-AnswerParagraphAfter($matchSubstring, $answerText) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like $matchSubstring) {
            $p.Range.InsertParagraphAfter()
            $newP = $d.Paragraphs.Item($i + 1)
            $newP.Range.Text = $answerText
            return
        }
    }
}

# --- Block 1: "Cau 1:" (Layer basics) true/false list ---
Set-ParagraphText "*Mỗi lớp ảnh chứa các đối tượng khác nhau*" `
    "a) Mỗi lớp ảnh chứa các đối tượng khác nhau và có thể được xử lý riêng biệt." | Out-Null

Set-ParagraphText "*Thứ tự sắp xếp các lớp không ảnh hưởng*" `
    "b) Thứ tự sắp xếp các lớp không ảnh hưởng đến kết quả cuối cùng của ảnh sản phẩm." | Out-Null

Set-ParagraphText "*Tạo một lớp mới trong GIMP*" `
    ("c) Tạo một lớp mới trong GIMP có thể thực hiện bằng lệnh " + $Q + "New Layer" + $Q + ".") | Out-Null

Set-ParagraphText "*Việc chỉnh sửa một lớp cụ thể trong GIMP*" `
    "d) Việc chỉnh sửa một lớp cụ thể trong GIMP yêu cầu người dùng chọn lớp đó trước khi thực hiện thao tác." | Out-Null

Insert-AnswerParagraphAfter "*Việc chỉnh sửa một lớp cụ thể trong GIMP*" "Đáp án: ĐSĐĐ"

# --- Block 2: "Cau 2:" (Layer commands) true/false list ---
Set-ParagraphText "*Tạo một bản sao của lớp được chọn*" `
    ("a) Tạo một bản sao của lớp được chọn có thể thực hiện bằng lệnh " + $Q + "Duplicate Layer" + $Q + ".") | Out-Null

Set-ParagraphText "*được sử dụng để xóa toàn bộ ảnh*" `
    ("b) Lệnh " + $Q + "Delete Layer" + $Q + " được sử dụng để xóa toàn bộ ảnh, không chỉ lớp được chọn.") | Out-Null

Set-ParagraphText "*là lệnh dùng để gộp lớp hiện tại*" `
    ("c) " + $Q + "Merge Down" + $Q + " là lệnh dùng để gộp lớp hiện tại với lớp ngay dưới nó.") | Out-Null

Set-ParagraphText "*Xóa một lớp ảnh trong GIMP có thể thực hiện mà không ảnh hưởng*" `
    "d) Xóa một lớp ảnh trong GIMP có thể thực hiện mà không ảnh hưởng đến các lớp khác." | Out-Null

Insert-AnswerParagraphAfter "*Xóa một lớp ảnh trong GIMP có thể thực hiện mà không ảnh hưởng*" "Đáp án: ĐSĐĐ"
